$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025938051212694
$ws.Range("D2").Value = 1.028769038732889
$ws.Range("E2").Value = 1.026200958928118
$ws.Range("I2").Value = 1.028984268334685
$ws.Range("J2").Value = 1.031104281575443
$ws.Range("K2").Value = 1.031584893355141
$ws.Range("L2").Value = 1.029024298664871
$ws.Range("N2").Value = 1.032568567924598
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027269169358222
$ws.Range("D3").Value = 1.029710552217391
$ws.Range("E3").Value = 1.027342097813834
$ws.Range("I3").Value = 1.02920245316513
$ws.Range("J3").Value = 1.032072980556356
$ws.Range("K3").Value = 1.0323340519256
$ws.Range("L3").Value = 1.029971997544187
$ws.Range("N3").Value = 1.033538642569172
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028129570331277
$ws.Range("D4").Value = 1.030318776568588
$ws.Range("E4").Value = 1.02808001213684
$ws.Range("I4").Value = 1.029341852899161
$ws.Range("J4").Value = 1.032698455911838
$ws.Range("K4").Value = 1.032817193912261
$ws.Range("L4").Value = 1.030584172063448
$ws.Range("N4").Value = 1.034165006171401
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028491066461249
$ws.Range("D5").Value = 1.030574237052479
$ws.Range("E5").Value = 1.028390119350224
$ws.Range("I5").Value = 1.029400030896381
$ws.Range("J5").Value = 1.032961088131471
$ws.Range("K5").Value = 1.033019922299587
$ws.Range("L5").Value = 1.030841281284207
$ws.Range("N5").Value = 1.03442801135891
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028551750680734
$ws.Range("D6").Value = 1.030617116145508
$ws.Range("E6").Value = 1.02844218122728
$ws.Range("I6").Value = 1.029409774294766
$ws.Range("J6").Value = 1.033005166687222
$ws.Range("K6").Value = 1.033053938793
$ws.Range("L6").Value = 1.030884436525717
$ws.Range("N6").Value = 1.034472152511267
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028134401505556
$ws.Range("D7").Value = 1.030322190973187
$ws.Range("E7").Value = 1.028084156241262
$ws.Range("I7").Value = 1.029342631948205
$ws.Range("J7").Value = 1.032701966461211
$ws.Range("K7").Value = 1.0328199042875
$ws.Range("L7").Value = 1.030587608544632
$ws.Range("N7").Value = 1.034168521706157
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026388102131658
$ws.Range("D8").Value = 1.029087435348456
$ws.Range("E8").Value = 1.02658671291504
$ws.Range("I8").Value = 1.029058373928939
$ws.Range("J8").Value = 1.031431936483775
$ws.Range("K8").Value = 1.031838410229246
$ws.Range("L8").Value = 1.029344796820987
$ws.Range("N8").Value = 1.032896688140478
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023303654946129
$ws.Range("D9").Value = 1.026903912191093
$ws.Range("E9").Value = 1.023944238815288
$ws.Range("I9").Value = 1.028543812040642
$ws.Range("J9").Value = 1.02918361843788
$ws.Range("K9").Value = 1.030096450915342
$ws.Range("L9").Value = 1.027146652176641
$ws.Range("N9").Value = 1.030645177225073
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021242221962613
$ws.Range("D10").Value = 1.025442915384209
$ws.Range("E10").Value = 1.022179871956603
$ws.Range("I10").Value = 1.028191545588298
$ws.Range("J10").Value = 1.027677604631381
$ws.Range("K10").Value = 1.028926662140312
$ws.Range("L10").Value = 1.025675585978313
$ws.Range("N10").Value = 1.029137024706227
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020348321341924
$ws.Range("D11").Value = 1.02480899938972
$ws.Range("E11").Value = 1.021415198240482
$ws.Range("I11").Value = 1.028036813104837
$ws.Range("J11").Value = 1.027023753587938
$ws.Range("K11").Value = 1.028418091757648
$ws.Range("L11").Value = 1.025037225774362
$ws.Range("N11").Value = 1.028482245119304
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020016088344602
$ws.Range("D12").Value = 1.024573337750742
$ws.Range("E12").Value = 1.021131057082707
$ws.Range("I12").Value = 1.027979007375359
$ws.Range("J12").Value = 1.026780619830912
$ws.Range("K12").Value = 1.028228876540236
$ws.Range("L12").Value = 1.024799900071498
$ws.Range("N12").Value = 1.028238766084455
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020087362517248
$ws.Range("D13").Value = 1.024623896927839
$ws.Range("E13").Value = 1.021192011236364
$ws.Range("I13").Value = 1.027991421901645
$ws.Range("J13").Value = 1.02683278486486
$ws.Range("K13").Value = 1.028269477892885
$ws.Range("L13").Value = 1.024850816807664
$ws.Range("N13").Value = 1.028291005198735
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020320862938479
$ws.Range("D14").Value = 1.02478952355865
$ws.Range("E14").Value = 1.021391713252485
$ws.Range("I14").Value = 1.028032041626461
$ws.Range("J14").Value = 1.027003661488393
$ws.Range("K14").Value = 1.028402457497243
$ws.Range("L14").Value = 1.025017612666066
$ws.Range("N14").Value = 1.028462124486672
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020464703793899
$ws.Range("D15").Value = 1.024891545411063
$ws.Range("E15").Value = 1.021514741934872
$ws.Range("I15").Value = 1.02805702486724
$ws.Range("J15").Value = 1.027108909083015
$ws.Range("K15").Value = 1.028484349531275
$ws.Range("L15").Value = 1.025120353138078
$ws.Range("N15").Value = 1.028567521544951
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021301519652272
$ws.Range("D16").Value = 1.025484958810049
$ws.Range("E16").Value = 1.022230605976066
$ws.Range("I16").Value = 1.028201768275129
$ws.Range("J16").Value = 1.027720961684243
$ws.Range("K16").Value = 1.028960371007667
$ws.Range("L16").Value = 1.025717922524221
$ws.Range("N16").Value = 1.029180443331077
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021826084014763
$ws.Range("D17").Value = 1.025856842965766
$ws.Range("E17").Value = 1.022679461091937
$ws.Range("I17").Value = 1.028291972713831
$ws.Range("J17").Value = 1.028104418651253
$ws.Range("K17").Value = 1.029258417726303
$ws.Range("L17").Value = 1.026092390645512
$ws.Range("N17").Value = 1.029564444850963
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022131929755492
$ws.Range("D18").Value = 1.026073631949459
$ws.Range("E18").Value = 1.022941204228903
$ws.Range("I18").Value = 1.028344375344415
$ws.Range("J18").Value = 1.028327915235643
$ws.Range("K18").Value = 1.02943206627228
$ws.Range("L18").Value = 1.026310678742396
$ws.Range("N18").Value = 1.029788258826138
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022236194411597
$ws.Range("D19").Value = 1.026147530305067
$ws.Range("E19").Value = 1.023030440705521
$ws.Range("I19").Value = 1.028362207338202
$ws.Range("J19").Value = 1.028404093506835
$ws.Range("K19").Value = 1.029491242586529
$ws.Range("L19").Value = 1.026385086907792
$ws.Range("N19").Value = 1.029864545279213
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021769816065212
$ws.Range("D20").Value = 1.025816956238144
$ws.Range("E20").Value = 1.022631310086121
$ws.Range("I20").Value = 1.028282316575519
$ws.Range("J20").Value = 1.028063294685892
$ws.Range("K20").Value = 1.029226460535655
$ws.Range("L20").Value = 1.026052227522355
$ws.Range("N20").Value = 1.029523262484854
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020252108406069
$ws.Range("D21").Value = 1.024740756102692
$ws.Range("E21").Value = 1.021332908975759
$ws.Range("I21").Value = 1.028020089279436
$ws.Range("J21").Value = 1.026953349906856
$ws.Range("K21").Value = 1.02836330688241
$ws.Range("L21").Value = 1.024968501273737
$ws.Range("N21").Value = 1.028411741456919
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019296714263076
$ws.Range("D22").Value = 1.024062965252417
$ws.Range("E22").Value = 1.020515930179734
$ws.Range("I22").Value = 1.027853299946266
$ws.Range("J22").Value = 1.026253952333501
$ws.Range("K22").Value = 1.027818815842234
$ws.Range("L22").Value = 1.024285901262364
$ws.Range("N22").Value = 1.027711350658788
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019803297635392
$ws.Range("D23").Value = 1.024422384008476
$ws.Range("E23").Value = 1.020949086257157
$ws.Range("I23").Value = 1.027941900088114
$ws.Range("J23").Value = 1.026624862498879
$ws.Range("K23").Value = 1.028107631528156
$ws.Range("L23").Value = 1.024647876980591
$ws.Range("N23").Value = 1.028082787559145
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021795241511738
$ws.Range("D24").Value = 1.025834979718856
$ws.Range("E24").Value = 1.022653067658223
$ws.Range("I24").Value = 1.02828668042409
$ws.Range("J24").Value = 1.028081877352685
$ws.Range("K24").Value = 1.029240901223326
$ws.Range("L24").Value = 1.02607037591917
$ws.Range("N24").Value = 1.029541871541166
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024101944086345
$ws.Range("D25").Value = 1.027469333133986
$ws.Range("E25").Value = 1.024627848078591
$ws.Range("I25").Value = 1.028678461738373
$ws.Range("J25").Value = 1.029766107791799
$ws.Range("K25").Value = 1.030548275569554
$ws.Range("L25").Value = 1.029971997544187
$ws.Range("N25").Value = 1.031228493780687
